$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1645569620253164
$ws.Range("C2").Value = 0.620253164556962
$ws.Range("J2").Value = 0.01265822784810127
$ws.Range("P2").Value = 0.09810126582278481
$ws.Range("S2").Value = 0.1044303797468354
$ws.Range("C3").Value = 0.02884615384615385
$ws.Range("J3").Value = 0.01923076923076923
$ws.Range("P3").Value = 0.7884615384615384
$ws.Range("S3").Value = 0.1634615384615385
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.6111111111111112
$ws.Range("S4").Value = 0.3148148148148148
$ws.Range("B6").Value = 0.08502024291497975
$ws.Range("F6").Value = 0.05668016194331984
$ws.Range("J6").Value = 0.2591093117408907
$ws.Range("O6").Value = 0.02024291497975709
$ws.Range("Q6").Value = 0.1862348178137652
$ws.Range("R6").Value = 0.0728744939271255
$ws.Range("S6").Value = 0.319838056680162
$ws.Range("B7").Value = 0.1172839506172839
$ws.Range("D7").Value = 0.0308641975308642
$ws.Range("F7").Value = 0.06172839506172839
$ws.Range("J7").Value = 0.1049382716049383
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.07407407407407407
$ws.Range("S7").Value = 0.4444444444444444
$ws.Range("B8").Value = 0.102880658436214
$ws.Range("D8").Value = 0.02880658436213992
$ws.Range("F8").Value = 0.07407407407407407
$ws.Range("J8").Value = 0.08436213991769548
$ws.Range("O8").Value = 0.02469135802469136
$ws.Range("Q8").Value = 0.1707818930041152
$ws.Range("R8").Value = 0.102880658436214
$ws.Range("S8").Value = 0.411522633744856
$ws.Range("B9").Value = 0.1020408163265306
$ws.Range("D9").Value = 0.02040816326530612
$ws.Range("F9").Value = 0.06122448979591837
$ws.Range("J9").Value = 0.1122448979591837
$ws.Range("O9").Value = 0.01020408163265306
$ws.Range("Q9").Value = 0.2091836734693878
$ws.Range("R9").Value = 0.09693877551020408
$ws.Range("S9").Value = 0.3877551020408163
$ws.Range("B10").Value = 0.1153546375681995
$ws.Range("D10").Value = 0.02494154325798909
$ws.Range("F10").Value = 0.0740452065471551
$ws.Range("J10").Value = 0.08417770849571317
$ws.Range("O10").Value = 0.01948558067030397
$ws.Range("Q10").Value = 0.2244738893219018
$ws.Range("R10").Value = 0.09586905689789556
$ws.Range("S10").Value = 0.3616523772408418
$ws.Range("G11").Value = 0.1022727272727273
$ws.Range("J11").Value = 0.143939393939394
$ws.Range("K11").Value = 0.1893939393939394
$ws.Range("L11").Value = 0.5568181818181818
$ws.Range("S11").Value = 0.007575757575757576
$ws.Range("G12").Value = 0.7733333333333333
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 0.006666666666666667
$ws.Range("S12").Value = 0.02
$ws.Range("G13").Value = 0.6842105263157895
$ws.Range("J13").Value = 0.3157894736842105
$ws.Range("F15").Value = 0.01626016260162602
$ws.Range("H15").Value = 0.1747967479674797
$ws.Range("I15").Value = 0.06910569105691057
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.1056910569105691
$ws.Range("M15").Value = 0.004065040650406504
$ws.Range("O15").Value = 0.06097560975609756
$ws.Range("S15").Value = 0.2357723577235772
$ws.Range("F16").Value = 0.009216589861751152
$ws.Range("H16").Value = 0.2304147465437788
$ws.Range("I16").Value = 0.05069124423963134
$ws.Range("J16").Value = 0.4285714285714285
$ws.Range("K16").Value = 0.04147465437788019
$ws.Range("M16").Value = 0.03225806451612903
$ws.Range("O16").Value = 0.06451612903225806
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.02070393374741201
$ws.Range("H17").Value = 0.1966873706004141
$ws.Range("I17").Value = 0.09316770186335403
$ws.Range("J17").Value = 0.443064182194617
$ws.Range("K17").Value = 0.08281573498964803
$ws.Range("M17").Value = 0.01656314699792961
$ws.Range("O17").Value = 0.06625258799171843
$ws.Range("S17").Value = 0.08074534161490683
$ws.Range("F18").Value = 0.0179372197309417
$ws.Range("H18").Value = 0.1883408071748879
$ws.Range("I18").Value = 0.1031390134529148
$ws.Range("J18").Value = 0.4573991031390134
$ws.Range("K18").Value = 0.07623318385650224
$ws.Range("M18").Value = 0.01345291479820628
$ws.Range("O18").Value = 0.04932735426008968
$ws.Range("S18").Value = 0.09417040358744394
$ws.Range("F19").Value = 0.02295081967213115
$ws.Range("H19").Value = 0.2180327868852459
$ws.Range("I19").Value = 0.08278688524590164
$ws.Range("J19").Value = 0.3819672131147541
$ws.Range("K19").Value = 0.09754098360655737
$ws.Range("M19").Value = 0.02131147540983606
$ws.Range("N19").Value = 0.000819672131147541
$ws.Range("O19").Value = 0.0860655737704918
$ws.Range("S19").Value = 0.08852459016393442
